$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seth Phebus (row 24) gets half credit, and his remark changes from the
# generic "Didn't attend the lab" to reflect that he did show the result.
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = "Didn't attend the lab but showed me the result."

# Scroll the view down a bit and move the active selection, matching
# where the instructor's cursor ended up after making the edit.
$excel.Goto($ws.Range("A3"), $true)
$ws.Range("H24").Select()
